# Habit-playbook-template.pptx edit
#
# Source diff changes the single "Ideas" content-placeholder shape
# (id 951 / creationId {F90B76D8-E2F5-4052-B20A-3341159EA7BF}) on the
# only slide so its text reads "Idea" instead of "Ideas".
#
# (The rest of the upstream diff -- the notesMaster "datetimeFigureOut"
# field text and the ppt/changesInfos/changesInfo1.xml collaboration
# audit trail -- are values PowerPoint recomputes/stamps automatically
# from the live system clock / co-authoring session and are not
# reachable through the PowerPoint object model, so there is nothing to
# script for those parts.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 951) {
        $targetShape = $sh
        break
    }
}

if ($targetShape -eq $null) {
    # Fallback: look the shape up by its placeholder name.
    $targetShape = $s.Shapes.Item("Content Placeholder 950")
}

$targetShape.TextFrame.TextRange.Text = "Idea"
